$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before column H (8) to hold "property_category" / "stock",
# shifting the existing date / legislator_name / legislator_id columns right.
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(2, 8).Value = "stock"
